# fix alignment tool output for benchmark accuracy
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: "Searching Time (s)" -> "Search Time (s)" ---
$ws.Range("D1").Value = "Search Time (s)"

# --- Fix mixed-up rows: row 5 (Boyer-Moore) and row 6 (Knuth-Morris-Pratt)
#     had swapped algorithm labels vs. their timing data. Swap the row
#     labels and the associated Search/Total time values back. ---
$ws.Range("A5").Value = "Knuth-Morris-Pratt algorithm"
$ws.Range("A6").Value = "Boyer-Moore algorithm"

$ws.Range("D5").Value = 2631.979
$ws.Range("E5").Value = 2631.979
$ws.Range("D6").Value = 3651.484
$ws.Range("E6").Value = 3651.484

# --- Mark the alignment-tool rows that are not exact/full search tools ---
$ws.Range("A9").Value = "* DIAMOND"
$ws.Range("A10").Value = "* Mmseqs2"

# --- Clean up stray formatted-but-empty rows (11-14) that extended all
#     the way out to column AMJ; clear them down to the end of the sheet ---
$lastRow = $ws.Rows.Count
$clearRange = $ws.Range("A11:G" + $lastRow)
$clearRange.EntireRow.Delete()

$ws.Range("A11").Select()
